# "tt call graph added"
#
# Adds two new slides at the end of the deck:
#  - a new slide (based on the existing trie-lookup call graph, slide 1)
#    redrawn/retexted into a typetable/trie hash-registration call graph
#  - a new slide that is an exact duplicate of the existing RPC call
#    graph (slide 2)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# New slide B: duplicate slide 2 (RPC call graph) unchanged. This is
# created first so it picks up the lower internal slide id (258),
# matching the id ordering seen in the target deck (ids are assigned
# in creation order, then the slides get reordered below).
# ---------------------------------------------------------------
$src2 = $p.Slides.Item(2)
$dupRange2 = $src2.Duplicate()
$slideB = $dupRange2.Item(1)

# ---------------------------------------------------------------
# New slide A: duplicate slide 1 (trie call graph), then reposition /
# retext its shapes and add one new shape.
# ---------------------------------------------------------------
$src1 = $p.Slides.Item(1)
$dupRange1 = $src1.Duplicate()
$slideA = $dupRange1.Item(1)
$slideA.MoveTo($p.Slides.Count)

# Shape "Rettangolo 3" (id 4): trie.Lookup() -> typetable.TypeForHash()
$sh = $slideA.Shapes.Item(1)
$sh.Left = 123.56590551181102
$sh.Top = 241.61
$sh.Width = 207.70622047244095
$sh.Height = 25.960314960629923
$sh.TextFrame.TextRange.Text = "typetable.TypeForHash()"

# Shape "Rettangolo 13" (id 14): trie.prefix() -> trie.Lookup()
$sh = $slideA.Shapes.Item(2)
$sh.Left = 167.79968503937008
$sh.Top = 375.513937007874
$sh.Width = 119.23866141732283
$sh.Height = 25.960314960629923
$sh.TextFrame.TextRange.Text = "trie.Lookup()"

# Connector "Connettore 2 15" (id 16): stCxn id4/idx2 -> endCxn id14/idx0 (unchanged wiring, new geometry)
$sh = $slideA.Shapes.Item(3)
$sh.Left = 227.41905511811024
$sh.Top = 267.5703149606299
$sh.Width = 0.0
$sh.Height = 107.94362204724409

# Shape "Rettangolo 18" (id 19): trie.Delete() -> newTypeTable()
$sh = $slideA.Shapes.Item(4)
$sh.Left = 383.4289763779528
$sh.Top = 122.84685039370079
$sh.Width = 158.20370078740157
$sh.Height = 25.960314960629923
$sh.TextFrame.TextRange.Text = "newTypeTable()"

# Shape "Rettangolo 21" (id 22): trie.Insert() -> typetable.Register()
$sh = $slideA.Shapes.Item(5)
$sh.Left = 374.54417322834644
$sh.Top = 242.16574803149607
$sh.Width = 175.97330708661417
$sh.Height = 25.960314960629923
$sh.TextFrame.TextRange.Text = "typetable.Register()"

# Shape "Rettangolo 22" (id 23): trie.prefix() -> trie.Insert()
$sh = $slideA.Shapes.Item(6)
$sh.Left = 402.9114960629921
$sh.Top = 374.81472440944884
$sh.Width = 119.23866141732283
$sh.Height = 25.960314960629923
$sh.TextFrame.TextRange.Text = "trie.Insert()"

# Connector "Connettore 2 23" (id 24): stCxn id22/idx2 -> endCxn id23/idx0 (unchanged wiring, new geometry)
$sh = $slideA.Shapes.Item(7)
$sh.Left = 462.5308661417323
$sh.Top = 268.126062992126
$sh.Width = 0.0
$sh.Height = 106.68866141732283

# New shape "Rettangolo 16" (id 17): typetable.HashForType() -- built by
# duplicating an existing styled rectangle so fill/line/style match.
$newRange = $slideA.Shapes.Item(1).Duplicate()
$newSh = $newRange.Item(1)
$newSh.Name = "Rettangolo 16"
$newSh.Left = 123.56590551181102
$newSh.Top = 123.47732283464568
$newSh.Width = 207.70622047244095
$newSh.Height = 25.960314960629923
$newSh.TextFrame.TextRange.Text = "typetable.HashForType()"

# Move slide B (RPC duplicate) to the very end, after slide A, so the
# final order is: slide1, slide2, slideA (typetable graph), slideB
# (RPC duplicate) -- matching the target slide id list (259, 258).
$slideB.MoveTo($p.Slides.Count)

Write-Host "Final slide count: $($p.Slides.Count)"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    Write-Host "Slide $i : Id=$($p.Slides.Item($i).SlideID) Shapes=$($p.Slides.Item($i).Shapes.Count)"
}
